$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 1.607001066207886
$ws.Range("E2").Value = 1973.890429266887
$ws.Range("F2").Value = 0.08345645117945431
$ws.Range("G2").Value = 0.06786062060373141
$ws.Range("H2").Value = 0.0614537570010119
$ws.Range("I2").Value = 0.05494461140583885
$ws.Range("J2").Value = 0.05115609017150262
$ws.Range("K2").Value = 0.04871408880806793
$ws.Range("L2").Value = 0.04644142098957423
$ws.Range("M2").Value = 0.04489241006248255
$ws.Range("N2").Value = 0.04257085065188287
$ws.Range("O2").Value = 0.04257085065188287
$ws.Range("P2").Value = 0.04152764992100125
$ws.Range("Q2").Value = 0.04084550187709777
$ws.Range("R2").Value = 0.04014512118814641
$ws.Range("S2").Value = 0.03991046290851418
$ws.Range("T2").Value = 0.03950128536558406
$ws.Range("U2").Value = 0.03920602850328609
$ws.Range("V2").Value = 0.03884451333721058
$ws.Range("W2").Value = 0.03864099715114228
$ws.Range("X2").Value = 0.03863285103262266
$ws.Range("Y2").Value = 0.03847739628200558
$ws.Range("C3").Value = 1.63303279876709
$ws.Range("E3").Value = 2014.515077347611
$ws.Range("F3").Value = 0.07836976113253634
$ws.Range("G3").Value = 0.06396942730775021
$ws.Range("H3").Value = 0.05713125856831189
$ws.Range("I3").Value = 0.05468053967761684
$ws.Range("J3").Value = 0.05094361862009548
$ws.Range("K3").Value = 0.04907778822216196
$ws.Range("L3").Value = 0.04712413122966248
$ws.Range("M3").Value = 0.04430535657609684
$ws.Range("N3").Value = 0.04351150647120578
$ws.Range("O3").Value = 0.04328628964204857
$ws.Range("P3").Value = 0.0420424679302205
$ws.Range("Q3").Value = 0.04136415263166243
$ws.Range("R3").Value = 0.04020081330373207
$ws.Range("S3").Value = 0.04020081330373207
$ws.Range("T3").Value = 0.03979251628099797
$ws.Range("U3").Value = 0.03968849689909708
$ws.Range("V3").Value = 0.03943644762289034
$ws.Range("W3").Value = 0.03933676917653871
$ws.Range("X3").Value = 0.03930644720450959
$ws.Range("Y3").Value = 0.03926929975336473
$ws.Range("C4").Value = 1.569000959396362
$ws.Range("E4").Value = 1999.060721186033
$ws.Range("F4").Value = 0.08238323697220018
$ws.Range("G4").Value = 0.06808018002366137
$ws.Range("H4").Value = 0.05892998329543316
$ws.Range("I4").Value = 0.05496894003486948
$ws.Range("J4").Value = 0.05204669372541512
$ws.Range("K4").Value = 0.0493630929780033
$ws.Range("L4").Value = 0.04672208311362328
$ws.Range("M4").Value = 0.04480723952004904
$ws.Range("N4").Value = 0.04346242808393879
$ws.Range("O4").Value = 0.04254142908880976
$ws.Range("P4").Value = 0.04169619462146426
$ws.Range("Q4").Value = 0.04104724378693023
$ws.Range("R4").Value = 0.04070207340586722
$ws.Range("S4").Value = 0.04013920904257353
$ws.Range("T4").Value = 0.03974295238924667
$ws.Range("U4").Value = 0.0395117429117357
$ws.Range("V4").Value = 0.03935161444010415
$ws.Range("W4").Value = 0.03909655044143595
$ws.Range("X4").Value = 0.03901713349519429
$ws.Range("Y4").Value = 0.03896804524729109
$ws.Range("C5").Value = 1.76601505279541
$ws.Range("E5").Value = 1957.336499657174
$ws.Range("F5").Value = 0.08332040970381824
$ws.Range("G5").Value = 0.06805263263081575
$ws.Range("H5").Value = 0.06042311185126605
$ws.Range("I5").Value = 0.0542972548434103
$ws.Range("J5").Value = 0.05052878850316003
$ws.Range("K5").Value = 0.04753851746812396
$ws.Range("L5").Value = 0.04584487131679202
$ws.Range("M5").Value = 0.04404503857394814
$ws.Range("N5").Value = 0.04251195801996977
$ws.Range("O5").Value = 0.04132068697631465
$ws.Range("P5").Value = 0.04132068697631465
$ws.Range("Q5").Value = 0.04008240528505153
$ws.Range("R5").Value = 0.03973252510869901
$ws.Range("S5").Value = 0.03930331442865177
$ws.Range("T5").Value = 0.03900467900429944
$ws.Range("U5").Value = 0.0387312541202845
$ws.Range("V5").Value = 0.0386406782294382
$ws.Range("W5").Value = 0.03840639619956998
$ws.Range("X5").Value = 0.03821959134011822
$ws.Range("Y5").Value = 0.03815470759565642
$ws.Range("C6").Value = 1.58466362953186
$ws.Range("E6").Value = 1907.604185457702
$ws.Range("F6").Value = 0.08190951544406262
$ws.Range("G6").Value = 0.06714228489507346
$ws.Range("H6").Value = 0.05801011039337731
$ws.Range("I6").Value = 0.05500119245494482
$ws.Range("J6").Value = 0.04904866766873454
$ws.Range("K6").Value = 0.04696334980099157
$ws.Range("L6").Value = 0.04499816284385403
$ws.Range("M6").Value = 0.04301047006232111
$ws.Range("N6").Value = 0.04148706175084179
$ws.Range("O6").Value = 0.04037169944723953
$ws.Range("P6").Value = 0.03973887293405204
$ws.Range("Q6").Value = 0.0390748271860047
$ws.Range("R6").Value = 0.03848682632382926
$ws.Range("S6").Value = 0.03823017856429317
$ws.Range("T6").Value = 0.03787184661291619
$ws.Range("U6").Value = 0.03774446855299713
$ws.Range("V6").Value = 0.0374625304218001
$ws.Range("W6").Value = 0.03732211898401153
$ws.Range("X6").Value = 0.03722718784262591
$ws.Range("Y6").Value = 0.0371852667730546
$ws.Range("C7").Value = 1.593999147415161
$ws.Range("E7").Value = 1979.806861188283
$ws.Range("F7").Value = 0.08256413882546727
$ws.Range("G7").Value = 0.06445404629527572
$ws.Range("H7").Value = 0.05824488330209657
$ws.Range("I7").Value = 0.05386955291848522
$ws.Range("J7").Value = 0.04868576967742164
$ws.Range("K7").Value = 0.04752597622224912
$ws.Range("L7").Value = 0.04524267591962723
$ws.Range("M7").Value = 0.04377314101677759
$ws.Range("N7").Value = 0.0428678976989168
$ws.Range("O7").Value = 0.04247406257943531
$ws.Range("P7").Value = 0.04109912148656076
$ws.Range("Q7").Value = 0.0407154382395223
$ws.Range("R7").Value = 0.04023811562785921
$ws.Range("S7").Value = 0.03990933955861808
$ws.Range("T7").Value = 0.03941630212269431
$ws.Range("U7").Value = 0.03934439423358857
$ws.Range("V7").Value = 0.03900405803973058
$ws.Range("W7").Value = 0.03880862967531601
$ws.Range("X7").Value = 0.03867189766348188
$ws.Range("Y7").Value = 0.03859272633895287
$ws.Range("C8").Value = 1.740999221801758
$ws.Range("E8").Value = 1949.22379304587
$ws.Range("F8").Value = 0.08067436146094809
$ws.Range("G8").Value = 0.06711024690728068
$ws.Range("H8").Value = 0.06046169731150988
$ws.Range("I8").Value = 0.05406550902210983
$ws.Range("J8").Value = 0.05042628068437714
$ws.Range("K8").Value = 0.04624426572072349
$ws.Range("L8").Value = 0.04476184784393068
$ws.Range("M8").Value = 0.04292203071439091
$ws.Range("N8").Value = 0.04177365503538249
$ws.Range("O8").Value = 0.04044989145849479
$ws.Range("P8").Value = 0.04025147629179446
$ws.Range("Q8").Value = 0.03976172952539592
$ws.Range("R8").Value = 0.03944991198266924
$ws.Range("S8").Value = 0.03901006154911706
$ws.Range("T8").Value = 0.03878513591460542
$ws.Range("U8").Value = 0.03845583347581545
$ws.Range("V8").Value = 0.03842366751751592
$ws.Range("W8").Value = 0.03815031277459
$ws.Range("X8").Value = 0.03806765605121677
$ws.Range("Y8").Value = 0.03799656516658616
$ws.Range("C9").Value = 1.726998805999756
$ws.Range("E9").Value = 1904.716559579532
$ws.Range("F9").Value = 0.07948649683433064
$ws.Range("G9").Value = 0.06715619517809811
$ws.Range("H9").Value = 0.05897172152490451
$ws.Range("I9").Value = 0.05058738102345852
$ws.Range("J9").Value = 0.04888273747299154
$ws.Range("K9").Value = 0.04635040676633503
$ws.Range("L9").Value = 0.04519827840495585
$ws.Range("M9").Value = 0.04315172542178676
$ws.Range("N9").Value = 0.04234512075079135
$ws.Range("O9").Value = 0.04067897585906567
$ws.Range("P9").Value = 0.03955632721906479
$ws.Range("Q9").Value = 0.03890572970457481
$ws.Range("R9").Value = 0.03876070864884971
$ws.Range("S9").Value = 0.03800016457913042
$ws.Range("T9").Value = 0.0378108013277523
$ws.Range("U9").Value = 0.03767031874011723
$ws.Range("V9").Value = 0.03745996182362595
$ws.Range("W9").Value = 0.0372360535791456
$ws.Range("X9").Value = 0.03714864181381117
$ws.Range("Y9").Value = 0.03712897776958152
$ws.Range("C10").Value = 1.577041387557983
$ws.Range("E10").Value = 1972.100557205593
$ws.Range("F10").Value = 0.08372771692623537
$ws.Range("G10").Value = 0.06758183360786713
$ws.Range("H10").Value = 0.05836315317037141
$ws.Range("I10").Value = 0.05496803551361742
$ws.Range("J10").Value = 0.05255110243240427
$ws.Range("K10").Value = 0.04747353973602018
$ws.Range("L10").Value = 0.04567865779216917
$ws.Range("M10").Value = 0.04455062161708696
$ws.Range("N10").Value = 0.04303131039093363
$ws.Range("O10").Value = 0.04212724951933181
$ws.Range("P10").Value = 0.04094492584932839
$ws.Range("Q10").Value = 0.04052744560634285
$ws.Range("R10").Value = 0.04003455020277479
$ws.Range("S10").Value = 0.03949785940811568
$ws.Range("T10").Value = 0.03932148357491735
$ws.Range("U10").Value = 0.03897042307001688
$ws.Range("V10").Value = 0.03871348596102207
$ws.Range("W10").Value = 0.03863329151853557
$ws.Range("X10").Value = 0.03849591014061121
$ws.Range("Y10").Value = 0.03844250598841311
$ws.Range("C11").Value = 1.520999431610107
$ws.Range("E11").Value = 2011.696667030641
$ws.Range("F11").Value = 0.08171356825910424
$ws.Range("G11").Value = 0.06342189551528188
$ws.Range("H11").Value = 0.05923350671507112
$ws.Range("I11").Value = 0.05471760018887832
$ws.Range("J11").Value = 0.05048803707205774
$ws.Range("K11").Value = 0.04747542440926444
$ws.Range("L11").Value = 0.04624937630345609
$ws.Range("M11").Value = 0.04461232773506477
$ws.Range("N11").Value = 0.043073971710505
$ws.Range("O11").Value = 0.04224704496665176
$ws.Range("P11").Value = 0.04161649295999655
$ws.Range("Q11").Value = 0.04139150959832033
$ws.Range("R11").Value = 0.04070138153601927
$ws.Range("S11").Value = 0.04059088259500328
$ws.Range("T11").Value = 0.04018597076290401
$ws.Range("U11").Value = 0.03986864766811057
$ws.Range("V11").Value = 0.0395967396055722
$ws.Range("W11").Value = 0.03936846783190044
$ws.Range("X11").Value = 0.03936846783190044
$ws.Range("Y11").Value = 0.0392143599811041
